# Applies the "MBA" staff-list update:
#  - six faculty records change title prefix "Mr./Mrs." -> "Dr."
#  - the last-used cell (A11) becomes the active selection, scrolled so
#    row 6 is the top visible row (matches topLeftCell="A6" in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MBA")

# --- Name column updates (column A, rows 2-11) ------------------------
$ws.Range("A2").Value  = "Dr. NIRMALRAJ A"
$ws.Range("A3").Value  = "Dr. KERINAB BEENU G H"
$ws.Range("A4").Value  = "Dr. REVATHI D"
$ws.Range("A5").Value  = "Dr. MUZHUMATHI R"
$ws.Range("A6").Value  = "Dr. NAVENA NESA KUMARI J"
$ws.Range("A9").Value  = "Dr. NAFEZA E "

# --- Selection / scroll position ---------------------------------------
$ws.Activate()
$ws.Range("A11").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
